$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update "Label" text for rows 2-3 (split "Best Recall, Best F1" -> "Best F1") ---
$ws.Range("T2").Value = "Best F1"
$ws.Range("T3").Value = "Best F1"

# --- Row 4 (Model=RNN_32->256_units_predictions, Algorithm=correction) ---
$ws.Range("A4").Value = "RNN_256_units_predictions"
$ws.Range("D4").Value = 0.1955307262569832
$ws.Range("E4").Value = 0.9859154929577464
$ws.Range("F4").Value = 0.3263403263403263
$ws.Range("G4").Value = 0.1955307262569832
$ws.Range("H4").Value = 0.9859154929577464
$ws.Range("I4").Value = 0.3263403263403263
$ws.Range("K4").Value = 0.1972222222222222
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.3294663573085847
$ws.Range("N4").Value = 0.001691495965239015
$ws.Range("O4").Value = 0.01408450704225361
$ws.Range("P4").Value = 0.00312603096825842
$ws.Range("Q4").Value = 0.00865079365079382
$ws.Range("R4").Value = 0.01428571428571438
$ws.Range("S4").Value = 0.009579052038449016
$ws.Range("T4").Value = "Best Recall"

# --- Row 5 (Model=RNN_32->256_units_predictions, Algorithm=detection_correction) ---
$ws.Range("A5").Value = "RNN_256_units_predictions"
$ws.Range("D5").Value = 0.1955307262569832
$ws.Range("E5").Value = 0.9859154929577464
$ws.Range("F5").Value = 0.3263403263403263
$ws.Range("G5").Value = 0.1955307262569832
$ws.Range("H5").Value = 0.9859154929577464
$ws.Range("I5").Value = 0.3263403263403263
$ws.Range("K5").Value = 0.1945205479452055
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.3256880733944954
$ws.Range("N5").Value = -0.001010178311777715
$ws.Range("O5").Value = 0.01408450704225361
$ws.Range("P5").Value = -0.0006522529458308779
$ws.Range("Q5").Value = -0.00516634050880603
$ws.Range("R5").Value = 0.01428571428571438
$ws.Range("S5").Value = -0.001998689384010333
$ws.Range("T5").Value = "Best Recall"

# --- Row 6 (Model=CNN_256->32_filters_5_kernels_predictions, Algorithm=correction) ---
$ws.Range("A6").Value = "CNN_32_filters_5_kernels_predictions"
$ws.Range("D6").Value = 0.8571428571428571
$ws.Range("E6").Value = 0.1690140845070422
$ws.Range("F6").Value = 0.2823529411764706
$ws.Range("G6").Value = 0.8571428571428571
$ws.Range("H6").Value = 0.1690140845070422
$ws.Range("I6").Value = 0.2823529411764706
$ws.Range("K6").Value = 0.8571428571428571
$ws.Range("L6").Value = 0.1690140845070423
$ws.Range("M6").Value = 0.2823529411764706
$ws.Range("R6").Value = [double]"3.284409781182756E-16"
$ws.Range("T6").Value = "Best Precision"

# --- Row 7 (Model=CNN_256->32_filters_5_kernels_predictions, Algorithm=detection_correction) ---
$ws.Range("A7").Value = "CNN_32_filters_5_kernels_predictions"
$ws.Range("D7").Value = 0.8571428571428571
$ws.Range("E7").Value = 0.1690140845070422
$ws.Range("F7").Value = 0.2823529411764706
$ws.Range("G7").Value = 0.8571428571428571
$ws.Range("H7").Value = 0.1690140845070422
$ws.Range("I7").Value = 0.2823529411764706
$ws.Range("K7").Value = 0.8571428571428571
$ws.Range("L7").Value = 0.1690140845070423
$ws.Range("M7").Value = 0.2823529411764706
$ws.Range("R7").Value = [double]"3.284409781182756E-16"
$ws.Range("T7").Value = "Best Precision"

# --- Row 8 (Model=CNN_Attention... -> LSTM_32_layers_predictions, Algorithm=correction) ---
$ws.Range("A8").Value = "LSTM_32_layers_predictions"
$ws.Range("D8").Value = 0.7142857142857143
$ws.Range("E8").Value = 0.0704225352112676
$ws.Range("F8").Value = 0.1282051282051282
$ws.Range("G8").Value = 0.7142857142857143
$ws.Range("H8").Value = 0.0704225352112676
$ws.Range("I8").Value = 0.1282051282051282
$ws.Range("K8").Value = 0.7142857142857143
$ws.Range("L8").Value = 0.07042253521126761
$ws.Range("M8").Value = 0.1282051282051282
$ws.Range("N8").Value = 0
$ws.Range("O8").Value = [double]"1.387778780781446E-17"
$ws.Range("P8").Value = [double]"2.775557561562891E-17"
$ws.Range("Q8").Value = 0
$ws.Range("R8").Value = [double]"1.970645868709653E-16"
$ws.Range("S8").Value = [double]"2.164934898019055E-16"
$ws.Range("T8").Value = "Worst F1"

# --- Row 9 (Model=CNN_Attention... -> LSTM_32_layers_predictions, Algorithm=detection_correction) ---
$ws.Range("A9").Value = "LSTM_32_layers_predictions"
$ws.Range("D9").Value = 0.7142857142857143
$ws.Range("E9").Value = 0.0704225352112676
$ws.Range("F9").Value = 0.1282051282051282
$ws.Range("G9").Value = 0.7142857142857143
$ws.Range("H9").Value = 0.0704225352112676
$ws.Range("I9").Value = 0.1282051282051282
$ws.Range("K9").Value = 0.7142857142857143
$ws.Range("L9").Value = 0.07042253521126761
$ws.Range("M9").Value = 0.1282051282051282
$ws.Range("O9").Value = [double]"1.387778780781446E-17"
$ws.Range("R9").Value = [double]"1.970645868709653E-16"
$ws.Range("S9").Value = [double]"2.164934898019055E-16"
$ws.Range("T9").Value = "Worst F1"
